$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 408.1111
$ws.Range("I31").Value = 384.14285
$ws.Range("J31").Value = 492
$ws.Range("K31").Value = 1152.42855
$ws.Range("L31").Value = 1476
$ws.Range("M31").Value = -922.4285500000001
$ws.Range("N31").Value = -1936

$ws.Range("H46").Value = 433839
$ws.Range("J46").Value = 500500
$ws.Range("L46").Value = 1501500
$ws.Range("N46").Value = -1501738

$ws.Range("H60").Value = 433839
$ws.Range("J60").Value = 500500
$ws.Range("L60").Value = 1501500
$ws.Range("N60").Value = -1502468

$ws.Range("H131").Value = 1212.174
$ws.Range("I131").Value = 687
$ws.Range("J131").Value = 1616.1538
$ws.Range("K131").Value = 2061
$ws.Range("L131").Value = 4848.4614
$ws.Range("M131").Value = 2979
$ws.Range("N131").Value = -14928.4614

$ws.Range("H132").Value = 805233.7
$ws.Range("I132").Value = 1826.9149
$ws.Range("J132").Value = 3502385
$ws.Range("K132").Value = 5480.7447
$ws.Range("L132").Value = 10507155
$ws.Range("M132").Value = -2950.7447
$ws.Range("N132").Value = -10512215

$ws.Range("H136").Value = 27640
$ws.Range("J136").Value = 27640
$ws.Range("L136").Value = 27640
$ws.Range("N136").Value = -37840

$ws.Range("H137").Value = 1819460.6
$ws.Range("I137").Value = 3126092.5
$ws.Range("J137").Value = 1537.9131
$ws.Range("K137").Value = 9378277.5
$ws.Range("L137").Value = 4613.7393
$ws.Range("M137").Value = -9375727.5
$ws.Range("N137").Value = -9713.739300000001

$ws.Range("H141").Value = 1276.8541
$ws.Range("I141").Value = 1301.909
$ws.Range("K141").Value = 3905.727
$ws.Range("M141").Value = 1274.273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 200401180
$ws.Range("I61").Value = 333667330
$ws.Range("J61").Value = 502000
$ws.Range("K61").Value = 333667330
$ws.Range("L61").Value = 502000
$ws.Range("M61").Value = -333667118
$ws.Range("N61").Value = -502424

$ws.Range("H74").Value = 5251414.5
$ws.Range("I74").Value = 10041486
$ws.Range("J74").Value = 44815.13
$ws.Range("K74").Value = 10041486
$ws.Range("L74").Value = 44815.13
$ws.Range("M74").Value = -10040612
$ws.Range("N74").Value = -46563.13

$ws.Range("H77").Value = 5251414.5
$ws.Range("I77").Value = 10041486
$ws.Range("J77").Value = 44815.13
$ws.Range("K77").Value = 50207430
$ws.Range("L77").Value = 224075.65
$ws.Range("M77").Value = -50203062
$ws.Range("N77").Value = -232811.65

$ws.Range("H132").Value = 175318.17
$ws.Range("I132").Value = 112313.11
$ws.Range("K132").Value = 336939.33
$ws.Range("M132").Value = -334409.33

$ws.Range("H136").Value = 200401180
$ws.Range("I136").Value = 333667330
$ws.Range("J136").Value = 502000
$ws.Range("K136").Value = 1001001990
$ws.Range("L136").Value = 1506000
$ws.Range("M136").Value = -1000999440
$ws.Range("N136").Value = -1511100

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3279.3142
$ws.Range("I31").Value = 1165.7778
$ws.Range("J31").Value = 10412.5
$ws.Range("K31").Value = 1165.7778
$ws.Range("L31").Value = 10412.5
$ws.Range("M31").Value = -870.7778000000001
$ws.Range("N31").Value = -11002.5

$ws.Range("H34").Value = 3279.3142
$ws.Range("I34").Value = 1165.7778
$ws.Range("J34").Value = 10412.5
$ws.Range("K34").Value = 1165.7778
$ws.Range("L34").Value = 10412.5
$ws.Range("M34").Value = -963.7778000000001
$ws.Range("N34").Value = -10816.5

$ws.Range("H50").Value = 23228
$ws.Range("J50").Value = 23228
$ws.Range("L50").Value = 23228
$ws.Range("N50").Value = -24478

$ws.Range("H58").Value = 18520016
$ws.Range("I58").Value = 20834570
$ws.Range("J58").Value = 3583.5
$ws.Range("K58").Value = 20834570
$ws.Range("L58").Value = 3583.5
$ws.Range("M58").Value = -20834367
$ws.Range("N58").Value = -3989.5

$ws.Range("H132").Value = 58661.914
$ws.Range("I132").Value = 36078.586
$ws.Range("J132").Value = 167814.67
$ws.Range("K132").Value = 108235.758
$ws.Range("L132").Value = 503444.01
$ws.Range("M132").Value = -105705.758
$ws.Range("N132").Value = -508504.01

$ws.Range("H136").Value = 18520016
$ws.Range("I136").Value = 20834570
$ws.Range("J136").Value = 3583.5
$ws.Range("K136").Value = 62503710
$ws.Range("L136").Value = 10750.5
$ws.Range("M136").Value = -62501160
$ws.Range("N136").Value = -15850.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 537
$ws.Range("I6").Value = 76.14286
$ws.Range("J6").Value = 2150
$ws.Range("K6").Value = 228.42858
$ws.Range("L6").Value = 6450
$ws.Range("M6").Value = -115.42858
$ws.Range("N6").Value = -6676

$ws.Range("H7").Value = 466
$ws.Range("I7").Value = 116.36364
$ws.Range("J7").Value = 893.3333
$ws.Range("K7").Value = 349.09092
$ws.Range("L7").Value = 2679.9999
$ws.Range("M7").Value = -237.09092
$ws.Range("N7").Value = -2903.9999

$ws.Range("H56").Value = 154931.05
$ws.Range("I56").Value = 154931.05
$ws.Range("K56").Value = 154931.05
$ws.Range("M56").Value = -154401.05

$ws.Range("H98").Value = 422.1111
$ws.Range("I98").Value = 116
$ws.Range("J98").Value = 575.1667
$ws.Range("K98").Value = 348
$ws.Range("L98").Value = 1725.5001
$ws.Range("M98").Value = 1150
$ws.Range("N98").Value = -4721.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 2777
$ws.Range("I31").Value = 2777
$ws.Range("K31").Value = 2777
$ws.Range("M31").Value = -2485

$ws.Range("H37").Value = 2777
$ws.Range("I37").Value = 2777
$ws.Range("K37").Value = 2777
$ws.Range("M37").Value = -2500

$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H132").Value = 97156.71000000001
$ws.Range("I132").Value = 60725.94
$ws.Range("K132").Value = 182177.82
$ws.Range("M132").Value = -179647.82

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3012.4688
$ws.Range("J122").Value = 3297.5
$ws.Range("L122").Value = 9892.5
$ws.Range("N122").Value = -14792.5

$ws.Range("H132").Value = 31804.426
$ws.Range("I132").Value = 21702.078
$ws.Range("K132").Value = 65106.234
$ws.Range("M132").Value = -62576.234

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 30086.666
$ws.Range("I75").Value = 10000
$ws.Range("J75").Value = 40130
$ws.Range("K75").Value = 10000
$ws.Range("L75").Value = 40130
$ws.Range("N75").Value = -42002
$ws.Range("M75").Value = -9064

$ws.Range("H78").Value = 30086.666
$ws.Range("I78").Value = 10000
$ws.Range("J78").Value = 40130
$ws.Range("K78").Value = 30000
$ws.Range("L78").Value = 120390
$ws.Range("N78").Value = -129750
$ws.Range("M78").Value = -25320

$ws.Range("H132").Value = 125169.4
$ws.Range("I132").Value = 101136.5
$ws.Range("J132").Value = 221301
$ws.Range("K132").Value = 303409.5
$ws.Range("L132").Value = 663903
$ws.Range("M132").Value = -300879.5
$ws.Range("N132").Value = -668963
